$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "Работа по созданию функционала загрузки данных"

# Copy the date style (numFmt, border) from an existing date cell (C4)
# onto C17:C19 so the new rows reuse the existing style index instead of
# Excel fabricating a brand-new numFmt/style entry.
$ws.Range("C4").Copy()
$ws.Range("C17:C19").PasteSpecial(-4122)

# Row 17
$ws.Range("A17").Value = $newText
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = "3/25/2019"

# Row 18
$ws.Range("A18").Value = $newText
$ws.Range("B18").Value = 2
$ws.Range("C18").Value = "3/26/2019"

# Row 19
$ws.Range("A19").Value = $newText
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "3/27/2019"

# Update the selected cell shown in the saved view
[void]$ws.Range("A25").Select()
